$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value changes (numeric / text) for rows 10-21 (species records permuted) ---
$ws.Range("A10").Value = 111936796
$ws.Range("B10").Value = 56398
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 448882.8980770012
$ws.Range("R10").Value = 7087229.443335658
$ws.Range("A11").Value = 111936866
$ws.Range("B11").Value = 89423
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma"
$ws.Range("H11").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q11").Value = 448765.5992023234
$ws.Range("R11").Value = 7087416.731054713
$ws.Range("A12").Value = 111936893
$ws.Range("B12").Value = 77515
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("Q12").Value = 448742.3011697636
$ws.Range("R12").Value = 7087501.648173723
$ws.Range("A13").Value = 111936870
$ws.Range("B13").Value = 89423
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = "Granticka"
$ws.Range("G13").Value = "Porodaedalea chrysoloma"
$ws.Range("H13").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 449019.027096529
$ws.Range("R13").Value = 7087276.979166135
$ws.Range("A14").Value = 111936795
$ws.Range("B14").Value = 56398
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = "Tretåig hackspett"
$ws.Range("G14").Value = "Picoides tridactylus"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("Q14").Value = 448749.3706757246
$ws.Range("R14").Value = 7087421.839990681
$ws.Range("A15").Value = 111936869
$ws.Range("B15").Value = 89423
$ws.Range("E15").Value = 5432
$ws.Range("F15").Value = "Granticka"
$ws.Range("G15").Value = "Porodaedalea chrysoloma"
$ws.Range("H15").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q15").Value = 449143.8568242944
$ws.Range("R15").Value = 7087117.752608996
$ws.Range("A16").Value = 111936798
$ws.Range("B16").Value = 56398
$ws.Range("E16").Value = 100109
$ws.Range("F16").Value = "Tretåig hackspett"
$ws.Range("G16").Value = "Picoides tridactylus"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("Q16").Value = 448923.1251473558
$ws.Range("R16").Value = 7087371.00725084
$ws.Range("A17").Value = 111936868
$ws.Range("B17").Value = 89423
$ws.Range("E17").Value = 5432
$ws.Range("F17").Value = "Granticka"
$ws.Range("G17").Value = "Porodaedalea chrysoloma"
$ws.Range("H17").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q17").Value = 448988.017639213
$ws.Range("R17").Value = 7087186.778340456
$ws.Range("A19").Value = 111936792
$ws.Range("B19").Value = 90087
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 3298
$ws.Range("F19").Value = "Trådticka"
$ws.Range("G19").Value = "Climacocystis borealis"
$ws.Range("H19").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q19").Value = 448761.1990147882
$ws.Range("R19").Value = 7087578.827763715
$ws.Range("A20").Value = 111936865
$ws.Range("Q20").Value = 448738.4239939091
$ws.Range("R20").Value = 7087426.42220111
$ws.Range("A21").Value = 111936867
$ws.Range("Q21").Value = 448791.554596175
$ws.Range("R21").Value = 7087501.648173723

# --- Cells that need to become present-but-empty (K/L/M/N) or carry text (AC) ---
$ws.Range("K10").Font.Bold = $false
$ws.Range("L10").Font.Bold = $false
$ws.Range("M10").Font.Bold = $false
$ws.Range("N10").Font.Bold = $false
$ws.Range("AC10").Value = "ringhack äldre"
$ws.Range("K14").Font.Bold = $false
$ws.Range("L14").Font.Bold = $false
$ws.Range("M14").Font.Bold = $false
$ws.Range("N14").Font.Bold = $false
$ws.Range("AC14").Value = "ringhack äldre"
$ws.Range("K16").Font.Bold = $false
$ws.Range("L16").Font.Bold = $false
$ws.Range("M16").Font.Bold = $false
$ws.Range("N16").Font.Bold = $false
$ws.Range("AC16").Value = "ringhack äldre"

# --- Cells that must be fully cleared (removed) ---
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("AC13").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("AC15").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("AC17").ClearContents()